$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.7321483333333333
$ws.Range("H2").Value = 2.196445
$ws.Range("I2").Value = 0.05113520435363902
$ws.Range("J2").Value = 0.05113520435363902
$ws.Range("M2").Value = 1.193633666666667
$ws.Range("N2").Value = 3.580901
$ws.Range("O2").Value = 0.0852504197348203
$ws.Range("P2").Value = 0.08525041973482028
$ws.Range("Q2").Value = 0.8739168996605555
$ws.Range("R2").Value = 7.865252096944999
$ws.Range("S2").Value = 0.004359297634373536
$ws.Range("T2").Value = 0.004359297634373535

$ws.Range("G3").Value = 0.7321483333333333
$ws.Range("H3").Value = 2.196445
$ws.Range("I3").Value = 0.05113520435363902
$ws.Range("J3").Value = 0.05113520435363902
$ws.Range("O3").Value = 0.6175422122064692
$ws.Range("P3").Value = 0.6175422122064691
$ws.Range("Q3").Value = 6.330532766638889
$ws.Range("R3").Value = 56.97479489974999
$ws.Range("S3").Value = 0.03157814721817612
$ws.Range("T3").Value = 0.03157814721817611

$ws.Range("G4").Value = 0.7321483333333333
$ws.Range("H4").Value = 2.196445
$ws.Range("I4").Value = 0.05113520435363902
$ws.Range("J4").Value = 0.05113520435363902
$ws.Range("M4").Value = 4.028899666666667
$ws.Range("N4").Value = 12.086699
$ws.Range("O4").Value = 0.2877477380576656
$ws.Range("P4").Value = 0.2877477380576655
$ws.Range("Q4").Value = 2.949752176117222
$ws.Range("R4").Value = 26.547769585055
$ws.Range("S4").Value = 0.01471403938787612
$ws.Range("T4").Value = 0.01471403938787612

$ws.Range("G5").Value = 0.7321483333333333
$ws.Range("H5").Value = 2.196445
$ws.Range("I5").Value = 0.05113520435363902
$ws.Range("J5").Value = 0.05113520435363902
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.132449
$ws.Range("N5").Value = 0.397347
$ws.Range("O5").Value = 0.009459630001044888
$ws.Range("P5").Value = 0.009459630001044887
$ws.Range("Q5").Value = 0.09697231460166666
$ws.Range("R5").Value = 0.8727508314149999
$ws.Range("S5").Value = 0.0004837201132132448
$ws.Range("T5").Value = 0.0004837201132132448

$ws.Range("I6").Value = 0.7165747117895102
$ws.Range("J6").Value = 0.7165747117895102
$ws.Range("M6").Value = 1.193633666666667
$ws.Range("N6").Value = 3.580901
$ws.Range("O6").Value = 0.0852504197348203
$ws.Range("P6").Value = 0.08525041973482028
$ws.Range("Q6").Value = 12.24648964285756
$ws.Range("R6").Value = 110.218406785718
$ws.Range("S6").Value = 0.06108829495141362
$ws.Range("T6").Value = 0.06108829495141361

$ws.Range("I7").Value = 0.7165747117895102
$ws.Range("J7").Value = 0.7165747117895102
$ws.Range("O7").Value = 0.6175422122064692
$ws.Range("P7").Value = 0.6175422122064691
$ws.Range("S7").Value = 0.4425151327297072
$ws.Range("T7").Value = 0.4425151327297071

$ws.Range("I8").Value = 0.7165747117895102
$ws.Range("J8").Value = 0.7165747117895102
$ws.Range("M8").Value = 4.028899666666667
$ws.Range("N8").Value = 12.086699
$ws.Range("O8").Value = 0.2877477380576656
$ws.Range("P8").Value = 0.2877477380576655
$ws.Range("Q8").Value = 41.33586327012023
$ws.Range("R8").Value = 372.022769431082
$ws.Range("S8").Value = 0.2061927524667552
$ws.Range("T8").Value = 0.2061927524667551

$ws.Range("I9").Value = 0.7165747117895102
$ws.Range("J9").Value = 0.7165747117895102
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.132449
$ws.Range("N9").Value = 0.397347
$ws.Range("O9").Value = 0.009459630001044888
$ws.Range("P9").Value = 0.009459630001044887
$ws.Range("Q9").Value = 1.358905459860667
$ws.Range("R9").Value = 12.230149138746
$ws.Range("S9").Value = 0.006778531641634144
$ws.Range("T9").Value = 0.006778531641634143

$ws.Range("G10").Value = 2.568000333333333
$ws.Range("H10").Value = 7.704001
$ws.Range("I10").Value = 0.1793560346266988
$ws.Range("J10").Value = 0.1793560346266988
$ws.Range("M10").Value = 1.193633666666667
$ws.Range("N10").Value = 3.580901
$ws.Range("O10").Value = 0.0852504197348203
$ws.Range("P10").Value = 0.08525041973482028
$ws.Range("Q10").Value = 3.065251653877889
$ws.Range("R10").Value = 27.587264884901
$ws.Range("S10").Value = 0.01529017723389903
$ws.Range("T10").Value = 0.01529017723389903

$ws.Range("G11").Value = 2.568000333333333
$ws.Range("H11").Value = 7.704001
$ws.Range("I11").Value = 0.1793560346266988
$ws.Range("J11").Value = 0.1793560346266988
$ws.Range("O11").Value = 0.6175422122064692
$ws.Range("P11").Value = 0.6175422122064691
$ws.Range("Q11").Value = 22.20425768217222
$ws.Range("R11").Value = 199.83831913955
$ws.Range("S11").Value = 0.1107599223959517
$ws.Range("T11").Value = 0.1107599223959516

$ws.Range("G12").Value = 2.568000333333333
$ws.Range("H12").Value = 7.704001
$ws.Range("I12").Value = 0.1793560346266988
$ws.Range("J12").Value = 0.1793560346266988
$ws.Range("M12").Value = 4.028899666666667
$ws.Range("N12").Value = 12.086699
$ws.Range("O12").Value = 0.2877477380576656
$ws.Range("P12").Value = 0.2877477380576655
$ws.Range("Q12").Value = 10.34621568696656
$ws.Range("R12").Value = 93.11594118269899
$ws.Range("S12").Value = 0.05160929327082491
$ws.Range("T12").Value = 0.0516092932708249

$ws.Range("G13").Value = 2.568000333333333
$ws.Range("H13").Value = 7.704001
$ws.Range("I13").Value = 0.1793560346266988
$ws.Range("J13").Value = 0.1793560346266988
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.132449
$ws.Range("N13").Value = 0.397347
$ws.Range("O13").Value = 0.009459630001044888
$ws.Range("P13").Value = 0.009459630001044887
$ws.Range("Q13").Value = 0.3401290761496667
$ws.Range("R13").Value = 3.061161685347
$ws.Range("S13").Value = 0.001696641726023165
$ws.Range("T13").Value = 0.001696641726023165

$ws.Range("G14").Value = 0.7579039999999999
$ws.Range("H14").Value = 2.273712
$ws.Range("I14").Value = 0.05293404923015203
$ws.Range("J14").Value = 0.05293404923015203
$ws.Range("M14").Value = 1.193633666666667
$ws.Range("N14").Value = 3.580901
$ws.Range("O14").Value = 0.0852504197348203
$ws.Range("P14").Value = 0.08525041973482028
$ws.Range("Q14").Value = 0.9046597305013333
$ws.Range("R14").Value = 8.141937574511999
$ws.Range("S14").Value = 0.004512649915134101
$ws.Range("T14").Value = 0.004512649915134101

$ws.Range("G15").Value = 0.7579039999999999
$ws.Range("H15").Value = 2.273712
$ws.Range("I15").Value = 0.05293404923015203
$ws.Range("J15").Value = 0.05293404923015203
$ws.Range("O15").Value = 0.6175422122064692
$ws.Range("P15").Value = 0.6175422122064691
$ws.Range("Q15").Value = 6.553229567733332
$ws.Range("R15").Value = 58.97906610959999
$ws.Range("S15").Value = 0.03268900986263423
$ws.Range("T15").Value = 0.03268900986263423

$ws.Range("G16").Value = 0.7579039999999999
$ws.Range("H16").Value = 2.273712
$ws.Range("I16").Value = 0.05293404923015203
$ws.Range("J16").Value = 0.05293404923015203
$ws.Range("M16").Value = 4.028899666666667
$ws.Range("N16").Value = 12.086699
$ws.Range("O16").Value = 0.2877477380576656
$ws.Range("P16").Value = 0.2877477380576655
$ws.Range("Q16").Value = 3.053519172965333
$ws.Range("R16").Value = 27.481672556688
$ws.Range("S16").Value = 0.01523165293220936
$ws.Range("T16").Value = 0.01523165293220936

$ws.Range("G17").Value = 0.7579039999999999
$ws.Range("H17").Value = 2.273712
$ws.Range("I17").Value = 0.05293404923015203
$ws.Range("J17").Value = 0.05293404923015203
$ws.Range("K17").Value = 1
$ws.Range("L17").Value = 0.3333333333333333
$ws.Range("M17").Value = 0.132449
$ws.Range("N17").Value = 0.397347
$ws.Range("O17").Value = 0.009459630001044888
$ws.Range("P17").Value = 0.009459630001044887
$ws.Range("Q17").Value = 0.100383626896
$ws.Range("R17").Value = 0.9034526420639999
$ws.Range("S17").Value = 0.0005007365201743332
$ws.Range("T17").Value = 0.0005007365201743331
